$wb = $excel.ActiveWorkbook

# 1. Wreath sheet: remove the two sample/test rows (rows 2 and 3)
$wsWreath = $wb.Worksheets.Item("Wreath")
$wsWreath.Rows("2:3").Delete()

# 2. Candles sheet: add new product row (headers + one data row)
$wsCandles = $wb.Worksheets.Item("Candles")
$wsCandles.Range("A1").Value = "ธูป&เทียน"
$wsCandles.Range("B1").Value = "รายละเอียด"
$wsCandles.Range("C1").Value = "pathรูปภาพ"
$wsCandles.Range("D1").Value = "ราคา"
$wsCandles.Range("A2").Value = "ธูปแพ็ค24ก้าน"
$wsCandles.Range("B2").Value = "ทดสอบบบบบบบบบบ"
$wsCandles.Range("C2").Value = "C:\Users\User\Downloads\SnackBox4.png"
$wsCandles.Range("D2").Value = "60"

$wsCandles.Columns("A:B").ColumnWidth = 19.375
$wsCandles.Columns("C").ColumnWidth = 44.375
$wsCandles.Columns("D").ColumnWidth = 13.75

# 3. Switch the active/selected tab from SnackBox to Wreath
$wsWreath.Select()
$wsWreath.Range("A2:F3").Select()
